$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2367223065250379
$ws.Range("C2").Value = 0.4779969650986343
$ws.Range("J2").Value = 0.01213960546282246
$ws.Range("P2").Value = 0.1836115326251897
$ws.Range("S2").Value = 0.08952959028831563
$ws.Range("B3").Value = 0.0121580547112462
$ws.Range("C3").Value = 0.0182370820668693
$ws.Range("J3").Value = 0.02127659574468085
$ws.Range("P3").Value = 0.7689969604863222
$ws.Range("S3").Value = 0.1793313069908815
$ws.Range("J4").Value = 0.01052631578947368
$ws.Range("P4").Value = 0.7263157894736842
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.08421052631578947
$ws.Range("D6").Value = 0.008421052631578947
$ws.Range("E6").Value = 0.002105263157894737
$ws.Range("F6").Value = 0.07368421052631578
$ws.Range("J6").Value = 0.1978947368421053
$ws.Range("O6").Value = 0.02105263157894737
$ws.Range("Q6").Value = 0.1936842105263158
$ws.Range("R6").Value = 0.0568421052631579
$ws.Range("S6").Value = 0.3621052631578948
$ws.Range("B7").Value = 0.123943661971831
$ws.Range("D7").Value = 0.02253521126760564
$ws.Range("E7").Value = 0.002816901408450704
$ws.Range("F7").Value = 0.02816901408450704
$ws.Range("J7").Value = 0.1605633802816901
$ws.Range("O7").Value = 0.03380281690140845
$ws.Range("Q7").Value = 0.1971830985915493
$ws.Range("R7").Value = 0.04507042253521127
$ws.Range("S7").Value = 0.3859154929577465
$ws.Range("B8").Value = 0.08879023307436182
$ws.Range("D8").Value = 0.02663706992230854
$ws.Range("F8").Value = 0.06659267480577137
$ws.Range("J8").Value = 0.1187569367369589
$ws.Range("O8").Value = 0.02552719200887902
$ws.Range("Q8").Value = 0.2097669256381798
$ws.Range("R8").Value = 0.0732519422863485
$ws.Range("S8").Value = 0.390677025527192
$ws.Range("B9").Value = 0.08788598574821853
$ws.Range("D9").Value = 0.009501187648456057
$ws.Range("F9").Value = 0.07838479809976247
$ws.Range("J9").Value = 0.1116389548693587
$ws.Range("O9").Value = 0.02612826603325416
$ws.Range("Q9").Value = 0.2114014251781473
$ws.Range("R9").Value = 0.07600950118764846
$ws.Range("S9").Value = 0.3990498812351544
$ws.Range("B10").Value = 0.1096507698084867
$ws.Range("D10").Value = 0.02253098009763425
$ws.Range("F10").Value = 0.07134810364250845
$ws.Range("J10").Value = 0.1182876455125798
$ws.Range("O10").Value = 0.02027788208787082
$ws.Range("Q10").Value = 0.2260608336462636
$ws.Range("R10").Value = 0.06721742395794217
$ws.Range("S10").Value = 0.3646263612467142
$ws.Range("G11").Value = 0.153713298791019
$ws.Range("J11").Value = 0.1036269430051813
$ws.Range("K11").Value = 0.2107081174438687
$ws.Range("L11").Value = 0.5181347150259067
$ws.Range("S11").Value = 0.01381692573402418
$ws.Range("G12").Value = 0.7540983606557377
$ws.Range("J12").Value = 0.1770491803278688
$ws.Range("K12").Value = 0.01967213114754099
$ws.Range("L12").Value = 0.01639344262295082
$ws.Range("S12").Value = 0.03278688524590164
$ws.Range("G13").Value = 0.6081081081081081
$ws.Range("J13").Value = 0.3243243243243243
$ws.Range("S13").Value = 0.06756756756756757
$ws.Range("F15").Value = 0.02526315789473684
$ws.Range("H15").Value = 0.1494736842105263
$ws.Range("I15").Value = 0.06947368421052631
$ws.Range("J15").Value = 0.3831578947368421
$ws.Range("K15").Value = 0.0568421052631579
$ws.Range("M15").Value = 0.008421052631578947
$ws.Range("O15").Value = 0.0568421052631579
$ws.Range("S15").Value = 0.2505263157894737
$ws.Range("F16").Value = 0.03037383177570093
$ws.Range("H16").Value = 0.1822429906542056
$ws.Range("I16").Value = 0.07476635514018691
$ws.Range("J16").Value = 0.4299065420560748
$ws.Range("K16").Value = 0.1144859813084112
$ws.Range("M16").Value = 0.007009345794392523
$ws.Range("N16").Value = 0.002336448598130841
$ws.Range("O16").Value = 0.0514018691588785
$ws.Range("S16").Value = 0.1074766355140187
$ws.Range("F17").Value = 0.01832208293153327
$ws.Range("H17").Value = 0.1677917068466731
$ws.Range("I17").Value = 0.08968177434908389
$ws.Range("J17").Value = 0.4455159112825458
$ws.Range("K17").Value = 0.07714561234329798
$ws.Range("M17").Value = 0.01350048216007715
$ws.Range("N17").Value = 0.001928640308582449
$ws.Range("O17").Value = 0.0626808100289296
$ws.Range("S17").Value = 0.1234329797492768
$ws.Range("F18").Value = 0.03125
$ws.Range("H18").Value = 0.203125
$ws.Range("I18").Value = 0.11875
$ws.Range("J18").Value = 0.375
$ws.Range("K18").Value = 0.08437500000000001
$ws.Range("M18").Value = 0.009375
$ws.Range("O18").Value = 0.065625
$ws.Range("S18").Value = 0.1125
$ws.Range("F19").Value = 0.01572327044025157
$ws.Range("H19").Value = 0.2032232704402516
$ws.Range("I19").Value = 0.08883647798742138
$ws.Range("J19").Value = 0.3856132075471698
$ws.Range("K19").Value = 0.1002358490566038
$ws.Range("M19").Value = 0.02083333333333333
$ws.Range("N19").Value = 0.001572327044025157
$ws.Range("O19").Value = 0.07114779874213836
$ws.Range("S19").Value = 0.112814465408805
